$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 and J1, matching the style used by the other headers (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New data columns I and J
$ws.Cells.Item(2, 9).Value = 6
$ws.Cells.Item(2, 10).Value = 7

$ws.Cells.Item(3, 9).Value = 5
$ws.Cells.Item(3, 10).Value = 6

$ws.Cells.Item(4, 9).Value = 3
$ws.Cells.Item(4, 10).Value = 6

$ws.Cells.Item(5, 9).Value = 5
$ws.Cells.Item(5, 10).Value = 8
